{"js": "const doc = context.document;\nconst body = doc.body;\n\n// ------------------------------------------------------------------\n// 1) Relocate the \"_GoBack\" bookmark out of the\n//    \"Presentaci\u00f3n del cuento / en clase.\" sentence -- it will be\n//    re-inserted below, inside \"repositorio\".\n// ------------------------------------------------------------------\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// ------------------------------------------------------------------\n// 2) \"Documento escrito: correo lfrincon@javerianacali.edu.co\"\n//    -> \"Documento escrito: s\u00fabalo en su repositorio de github\"\n// ------------------------------------------------------------------\nconst mailResults = body.search(\": correo lfrincon@javerianacali.edu.co\", { matchCase: true });\nmailResults.load(\"text\");\nawait context.sync();\n\nmailResults.items[0].insertText(\": s\u00fabalo en su repositorio de github\", \"Replace\");\nawait context.sync();\n\n// Re-anchor \"_GoBack\" between \"reposi\" and \"torio de\" (mid-word, where\n// the author's cursor ended up).\nconst repoResults = body.search(\"s\u00fabalo en su reposi\", { matchCase: true });\nrepoResults.load(\"text\");\nawait context.sync();\n\nconst bookmarkSpot = repoResults.items[0].getRange(\"End\");\nbookmarkSpot.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// ------------------------------------------------------------------\n// 3) Split the delivery-date paragraph and add the new \"Presentaci\u00f3n\"\n//    paragraph right after it.\n// ------------------------------------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet dateParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"lunes 16 de marzo\") !== -1) {\n    dateParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\ndateParagraph\n  .getRange()\n  .insertText(\"Documento escrito:28 de marzo del 2020. En el repositorio de github\", \"Replace\");\nawait context.sync();\n\ndateParagraph.insertParagraph(\n  \" Presentaci\u00f3n: por definir seg\u00fan sigan las condiciones acad\u00e9micas.  Les estar\u00e9 informando. \",\n  \"After\"\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# ------------------------------------------------------------------\n# 1) \"Documento escrito: correo lfrincon@javerianacali.edu.co\"\n#    -> \"Documento escrito: s\u00fabalo en su repositorio de github\"\n# ------------------------------------------------------------------\n$find = $d.Content\n$find.Find.ClearFormatting()\n$find.Find.Replacement.ClearFormatting()\n$find.Find.Execute(\n    \": correo lfrincon@javerianacali.edu.co\",\n    $false, $false, $false, $false, $false,\n    $true, 1, $false,\n    \": s\u00fabalo en su repositorio de github\",\n    2\n) | Out-Null\n\n# ------------------------------------------------------------------\n# 2) Relocate the \"_GoBack\" bookmark: it used to sit inside the\n#    \"Presentaci\u00f3n del cuento / en clase.\" sentence; now it marks the\n#    point inside \"repositorio\" (between \"reposi\" and \"torio de\").\n# ------------------------------------------------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n$locate = $d.Content\n$locate.Find.ClearFormatting()\n$locate.Find.Execute(\"s\u00fabalo en su reposi\") | Out-Null\n$pos = $locate.End\n$target = $d.Range($pos, $pos)\n$d.Bookmarks.Add(\"_GoBack\", $target) | Out-Null\n\n# ------------------------------------------------------------------\n# 3) Split the delivery-date paragraph and add the new \"Presentaci\u00f3n\"\n#    paragraph right after it.\n# ------------------------------------------------------------------\n$idx = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -match \"lunes 16 de marzo\") {\n        $idx = $i\n    }\n}\n\n$p = $d.Paragraphs($idx)\n$r = $p.Range\n$r.MoveEnd(1, -1) | Out-Null   # wdCharacter -- exclude the paragraph mark\n$r.Text = \"Documento escrito:28 de marzo del 2020. En el repositorio de github\"\n\n$p.Range.InsertParagraphAfter() | Out-Null\n\n$newP = $d.Paragraphs($idx + 1)\n$newR = $newP.Range\n$newR.MoveEnd(1, -1) | Out-Null\n$newR.Text = \" Presentaci\u00f3n: por definir seg\u00fan sigan las condiciones acad\u00e9micas.  Les estar\u00e9 informando. \"\n"}
